$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.389.62"
$ws.Range("E2").Value = "  +6.69%  "

# Row 3
$ws.Range("D3").Value = "3.003.71"
$ws.Range("E3").Value = "  +4.23%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.79"
$ws.Range("E5").Value = "  +3.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.01"
$ws.Range("E6").Value = "  +7.18%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.517"
$ws.Range("E8").Value = "  +2.58%  "

# Row 9
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "2.995.83"
$ws.Range("E9").Value = "  +3.96%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.99"
$ws.Range("E10").Value = "  +1.45%  "

# Row 11
$ws.Range("E11").Value = "  +5.46%  "

# Row 12
$ws.Range("E12").Value = "  +3.93%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000241"
$ws.Range("E13").Value = "  +3.59%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.07"
$ws.Range("E14").Value = "  +6.61%  "

# Row 15
$ws.Range("E15").Value = "  +0.82%  "

# Row 16
$ws.Range("D16").Value = "65.341.87"
$ws.Range("E16").Value = "  +6.62%  "

# Row 17
$ws.Range("D17").Value = "3.500.50"
$ws.Range("E17").Value = "  +4.45%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.94"
$ws.Range("E18").Value = "  +5.70%  "

# Row 19
$ws.Range("D19").Value = "3.014.46"
$ws.Range("E19").Value = "  +4.45%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "451.37"
$ws.Range("E20").Value = "  +4.66%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.71"
$ws.Range("E21").Value = "  +5.07%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.683"
$ws.Range("E22").Value = "  +4.36%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.33"
$ws.Range("E23").Value = "  +7.39%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.38"
$ws.Range("E24").Value = "  +2.69%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.50"
$ws.Range("E25").Value = "  +5.33%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").Value = "  +10.88%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.72"
$ws.Range("E27").Value = "  +7.13%  "

# Row 28
$ws.Range("E28").Value = "  -0.01%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.42"
$ws.Range("E29").Value = "  +17.79%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.83"
$ws.Range("E30").Value = "  +11.88%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0000106"
$ws.Range("E31").Value = "  +1.36%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.61"
$ws.Range("E32").Value = "  +4.44%  "

# Row 33
$ws.Range("E33").Value = "  +5.00%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.95"
$ws.Range("E34").Value = "  +5.86%  "

# Row 35
$ws.Range("E35").Value = "  -0.06%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.984"
$ws.Range("E36").Value = "  +2.73%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.78"
$ws.Range("E37").Value = "  +7.46%  "

# Row 38
$ws.Range("E38").Value = "  +9.06%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "46.54"
$ws.Range("E39").Value = "  +19.61%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.21"
$ws.Range("E40").Value = "  +0.76%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.93"
$ws.Range("E41").Value = "  +3.37%  "

# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.121"
$ws.Range("E42").Value = "  +6.54%  "

# Row 43
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.300"
$ws.Range("E43").Value = "  +12.42%  "

# Row 44
$ws.Range("E44").Value = "  +2.41%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "382.83"
$ws.Range("E45").Value = "  +11.93%  "

# Row 46
$ws.Range("D46").Value = "2.769.32"
$ws.Range("E46").Value = "  +2.33%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0351"
$ws.Range("E47").Value = "  +4.57%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.88"
$ws.Range("E48").Value = "  +1.44%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.21"
$ws.Range("E50").Value = "  +7.87%  "

# Row 51
$ws.Range("E51").Value = "  +2.73%  "
